$d = $word.ActiveDocument
$dash = [char]0x2013

# Step 1: Remove the " - Corbin Peever" portion from the title, leaving the
# trailing period in place, e.g.:
#   "Publishing Manager Job Description - Corbin Peever."
#     -> "Publishing Manager Job Description."
$search = "Publishing Manager Job Description " + $dash + " Corbin Peever."
$replace = "Publishing Manager Job Description."
$found = $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
if (-not $found) {
    throw "Could not find the title text to replace."
}
Write-Host "Replaced title text: $found"

# Step 2: The target markup splits the title into two runs with identical
# formatting - one run holding "Publishing Manager Job Description" and a
# second run holding just the trailing ".". Re-applying the (already true)
# Bold formatting to the final character's range forces Word to break the
# run at that boundary, producing two <w:r> elements with matching <w:rPr>
# without altering any visible formatting.
$titlePara = $d.Paragraphs(1).Range
$periodRange = $d.Range($titlePara.End - 2, $titlePara.End - 1)
Write-Host "Period range text: [$($periodRange.Text)]"
$periodRange.Bold = 0
$periodRange.Bold = 1

Write-Host "Done"
